$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the obsolete "Level-a3" row (originally row 4) ---
$ws.Rows.Item(4).Delete()

# --- Insert 3 new rows at the top of the data block (before row 3) for the new levels ---
$ws.Rows.Item(3).Resize(3).Insert()

# Fill the new rows: index columns A/B and the level name in C
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Level-c1-t0-1216"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Level-c1-t1-1216"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "Level-c1-t2-1216"

# --- Renumber the A/B index columns for the rest of the (shifted) data rows ---
for ($r = 6; $r -le 35; $r++) {
  $n = $r - 2
  $ws.Cells.Item($r, 1).Value = $n
  $ws.Cells.Item($r, 2).Value = $n
}

# --- Match the existing Name-column formatting (style index 3) on the new C cells ---
$ws.Range("C6").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Give every index cell in column A (rows 3-35) the vertical-center style ---
$ws.Range("A3:A35").VerticalAlignment = -4108

# --- Re-anchor the duplicate-values conditional format to the new data range ---
$cf = $ws.Range("C3:C33").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("C6:C35"))
$cf.Priority = 2

# --- Update the remembered selection to match the new layout ---
$ws.Range("B3:B35").Select()

Write-Output "done"
